$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1490.9333
$ws.Range("I98").Value = 1356.7273
$ws.Range("K98").Value = 1356.7273
$ws.Range("M98").Value = 141.2727
$ws.Range("H122").Value = 1490.9333
$ws.Range("I122").Value = 1356.7273
$ws.Range("K122").Value = 4070.1819
$ws.Range("M122").Value = -1620.1819
$ws.Range("H137").Value = 3510.95
$ws.Range("I137").Value = 3828.8333
$ws.Range("J137").Value = 650
$ws.Range("K137").Value = 11486.4999
$ws.Range("L137").Value = 1950
$ws.Range("M137").Value = -8936.499899999999
$ws.Range("N137").Value = -7050
$ws.Range("H138").Value = 134698.47
$ws.Range("I138").Value = 2150.45
$ws.Range("J138").Value = 177455.9
$ws.Range("K138").Value = 6451.349999999999
$ws.Range("L138").Value = 532367.7
$ws.Range("M138").Value = -1311.349999999999
$ws.Range("N138").Value = -542647.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 470930.25
$ws.Range("I32").Value = 581275.8
$ws.Range("J32").Value = 12571.692
$ws.Range("K32").Value = 581275.8
$ws.Range("L32").Value = 12571.692
$ws.Range("M32").Value = -580988.8
$ws.Range("N32").Value = -13145.692
$ws.Range("H86").Value = 50028500
$ws.Range("J86").Value = 50028500
$ws.Range("L86").Value = 50028500
$ws.Range("N86").Value = -50030872
$ws.Range("H89").Value = 50028500
$ws.Range("J89").Value = 50028500
$ws.Range("L89").Value = 150085500
$ws.Range("N89").Value = -150097356
$ws.Range("H122").Value = 2934
$ws.Range("I122").Value = 1974.6666
$ws.Range("J122").Value = 3653.5
$ws.Range("K122").Value = 5923.9998
$ws.Range("L122").Value = 10960.5
$ws.Range("M122").Value = -3473.9998
$ws.Range("N122").Value = -15860.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 37436.75
$ws.Range("J93").Value = 37436.75
$ws.Range("L93").Value = 37436.75
$ws.Range("N93").Value = -41180.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1262725.5
$ws.Range("I6").Value = 1683466.6
$ws.Range("K6").Value = 1683466.6
$ws.Range("M6").Value = -1683353.6
$ws.Range("H31").Value = 2514.8125
$ws.Range("I31").Value = 1054.0741
$ws.Range("J31").Value = 10402.8
$ws.Range("K31").Value = 1054.0741
$ws.Range("L31").Value = 10402.8
$ws.Range("M31").Value = -759.0741
$ws.Range("N31").Value = -10992.8
$ws.Range("H34").Value = 2514.8125
$ws.Range("I34").Value = 1054.0741
$ws.Range("J34").Value = 10402.8
$ws.Range("K34").Value = 1054.0741
$ws.Range("L34").Value = 10402.8
$ws.Range("M34").Value = -852.0741
$ws.Range("N34").Value = -10806.8
$ws.Range("H132").Value = 15153848
$ws.Range("I132").Value = 972.4
$ws.Range("J132").Value = 27781244
$ws.Range("K132").Value = 2917.2
$ws.Range("L132").Value = 83343732
$ws.Range("M132").Value = -387.1999999999998
$ws.Range("N132").Value = -83348792
$ws.Range("H134").Value = 1643.4814
$ws.Range("I134").Value = 1471.5454
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 4414.6362
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -1879.6362
$ws.Range("N134").Value = -12270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 412.07407
$ws.Range("I107").Value = 380.4737
$ws.Range("J107").Value = 487.125
$ws.Range("K107").Value = 1141.4211
$ws.Range("L107").Value = 1461.375
$ws.Range("M107").Value = 778.5789
$ws.Range("N107").Value = -5301.375
$ws.Range("H122").Value = 6844.8125
$ws.Range("I122").Value = 376.58334
$ws.Range("J122").Value = 26249.5
$ws.Range("K122").Value = 3389.25006
$ws.Range("L122").Value = 236245.5
$ws.Range("M122").Value = -939.2500600000003
$ws.Range("N122").Value = -241145.5
$ws.Range("H132").Value = 3232.8816
$ws.Range("I132").Value = 2402.963
$ws.Range("J132").Value = 3690.1836
$ws.Range("K132").Value = 21626.667
$ws.Range("L132").Value = 33211.6524
$ws.Range("M132").Value = -19096.667
$ws.Range("N132").Value = -38271.6524
$ws.Range("H137").Value = 16320.625
$ws.Range("I137").Value = 20416.5
$ws.Range("J137").Value = 4033
$ws.Range("K137").Value = 61249.5
$ws.Range("L137").Value = 12099
$ws.Range("M137").Value = -56149.5
$ws.Range("N137").Value = -22299
$ws.Range("H139").Value = 4173.5454
$ws.Range("I139").Value = 1323.2222
$ws.Range("J139").Value = 17000
$ws.Range("K139").Value = 3969.6666
$ws.Range("L139").Value = 51000
$ws.Range("M139").Value = 1170.3334
$ws.Range("N139").Value = -61280
$ws.Range("H140").Value = 1940.7894
$ws.Range("J140").Value = 3250
$ws.Range("L140").Value = 9750
$ws.Range("N140").Value = -20110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41872
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -129360
$ws.Range("H98").Value = 61643
$ws.Range("J98").Value = 61643
$ws.Range("L98").Value = 61643
$ws.Range("N98").Value = -67633
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 35154.5
$ws.Range("J17").Value = 70009
$ws.Range("L17").Value = 70009
$ws.Range("N17").Value = -70349
$ws.Range("H23").Value = 49503.5
$ws.Range("J23").Value = 49503.5
$ws.Range("L23").Value = 49503.5
$ws.Range("N23").Value = -49963.5
$ws.Range("H61").Value = 3385.3794
$ws.Range("I61").Value = 3272.1904
$ws.Range("J61").Value = 3682.5
$ws.Range("K61").Value = 3272.1904
$ws.Range("L61").Value = 3682.5
$ws.Range("M61").Value = -3070.1904
$ws.Range("N61").Value = -4086.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H98").Value = 95000
$ws.Range("J98").Value = 95000
$ws.Range("L98").Value = 95000
$ws.Range("N98").Value = -100990
$ws.Range("H113").Value = 3385.3794
$ws.Range("I113").Value = 3272.1904
$ws.Range("J113").Value = 3682.5
$ws.Range("K113").Value = 3272.1904
$ws.Range("L113").Value = 3682.5
$ws.Range("M113").Value = -1102.1904
$ws.Range("N113").Value = -8022.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 47295.168
$ws.Range("J69").Value = 47295.168
$ws.Range("L69").Value = 47295.168
$ws.Range("N69").Value = -48793.168
$ws.Range("H72").Value = 47295.168
$ws.Range("J72").Value = 47295.168
$ws.Range("L72").Value = 141885.504
$ws.Range("N72").Value = -149373.504
$ws.Range("H132").Value = 3878644
$ws.Range("I132").Value = 3270.8
$ws.Range("J132").Value = 7248533.5
$ws.Range("K132").Value = 9812.400000000001
$ws.Range("L132").Value = 21745600.5
$ws.Range("M132").Value = -7282.400000000001
$ws.Range("N132").Value = -21750660.5
$ws.Range("H136").Value = 3163.3215
$ws.Range("I136").Value = 2618.2
$ws.Range("K136").Value = 7854.599999999999
$ws.Range("M136").Value = -5304.599999999999
